$wb = $excel.ActiveWorkbook

# --- HEB sheet ("HEB", 2nd tab): the data table was originally offset by
# three blank rows (header + 24 data rows living in A4:T28 instead of
# A1:T25). Remove the three leading rows so everything shifts up to align
# with the other sheets (header in row 1, data in rows 2:25).
$wsHEB = $wb.Worksheets.Item("HEB")
$wsHEB.Rows("1:3").Delete()

# --- UPE sheet ("UPE", 4th tab): just move the lingering selection to B21.
$wsUPE = $wb.Worksheets.Item("UPE")
$wsUPE.Range("B21").Select()

# --- Make HEB the active tab, with A2 selected (this also clears the
# tabSelected flag that used to sit on the IPE sheet, and drops HEB's old
# topLeftCell scroll-freeze now that the data starts at row 1).
$wsHEB.Activate()
$wsHEB.Range("A2").Select()
